# Apply the RPAR_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer cell (A18)
#    from 2021-03-24 to 2021-03-25
#  - refresh the Weight (col D) and Percent Change (col E) figures for
#    rows 2-15 with the newly-modeled values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no UI means to edit locked cells), so
# temporarily lift protection, make the edits, then restore it.
$ws.Unprotect()

# --- Disclaimer text (A18): date roll 2021-03-24 -> 2021-03-25 ---
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."
$ws.Range("A18").Value = $disclaimer
# Re-run autofit so the embedded line break doesn't leave a stray
# explicit row-height override behind (matches the original row 18,
# which carries no custom height).
$ws.Rows.Item(18).AutoFit()

# --- Weight (D) / Percent Change (E) refresh, rows 2-15 ---
$ws.Range("D2").Value = 0.05465660320273161
$ws.Range("E2").Value = 0.004052785054119434

$ws.Range("D3").Value = 0.02296087891426966
$ws.Range("E3").Value = 0.006382540662960734

$ws.Range("D4").Value = 0.03073066366995872
$ws.Range("E4").Value = 0.001582278481012667

$ws.Range("D5").Value = 0.03168086681217989
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.03231906295247769
$ws.Range("E6").Value = 0.002507522567702924

$ws.Range("D7").Value = 0.01846041428046593
$ws.Range("E7").Value = 0.01049935979513461

$ws.Range("D8").Value = 0.004538283664339895
$ws.Range("E8").Value = 0.004910714285714546

$ws.Range("D9").Value = 0.006454898104726298
$ws.Range("E9").Value = 0.01025319104415146

$ws.Range("D10").Value = 0.06973559094936572
$ws.Range("E10").Value = -0.003486345148169812

$ws.Range("D11").Value = 0.06989767250880644
$ws.Range("E11").Value = -0.003478260869565153

$ws.Range("D12").Value = 0.1493635597432628
$ws.Range("E12").Value = -0.007885408377342085

$ws.Range("D13").Value = 0.393021443390314
$ws.Range("E13").Value = -0.001314521076154485

$ws.Range("D14").Value = 0.1161800618071013
$ws.Range("E14").Value = -0.004150390624999889

$ws.Range("E15").Value = -0.001882847448836378

# Restore the original protection state (objects/scenarios locked,
# row/column formatting allowed).
$ws.Protect($null, $true, $true, $true)
